$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.062756666666667
$ws.Range("H2").Value = 12.18827
$ws.Range("I2").Value = 0.6829811567947219
$ws.Range("J2").Value = 0.6829811567947218
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 25.15544366666667
$ws.Range("N2").Value = 75.466331
$ws.Range("O2").Value = 0.9701024243751556
$ws.Range("P2").Value = 0.9701024243751556
$ws.Range("Q2").Value = 102.2004464597078
$ws.Range("R2").Value = 919.80401813737
$ws.Range("S2").Value = 0.662561676009108
$ws.Range("T2").Value = 0.6625616760091079

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.062756666666667
$ws.Range("H3").Value = 12.18827
$ws.Range("I3").Value = 0.6829811567947219
$ws.Range("J3").Value = 0.6829811567947218
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6301496666666666
$ws.Range("N3").Value = 1.890449
$ws.Range("O3").Value = 0.02430128951224074
$ws.Range("P3").Value = 0.02430128951224074
$ws.Range("Q3").Value = 2.560144759247778
$ws.Range("R3").Value = 23.04130283323
$ws.Range("S3").Value = 0.01659732282267363
$ws.Range("T3").Value = 0.01659732282267362

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.062756666666667
$ws.Range("H4").Value = 12.18827
$ws.Range("I4").Value = 0.6829811567947219
$ws.Range("J4").Value = 0.6829811567947218
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1451156666666667
$ws.Range("N4").Value = 0.435347
$ws.Range("O4").Value = 0.005596286112603657
$ws.Range("P4").Value = 0.005596286112603657
$ws.Range("Q4").Value = 0.5895696421877779
$ws.Range("R4").Value = 5.30612677969
$ws.Range("S4").Value = 0.003822157962940283
$ws.Range("T4").Value = 0.003822157962940283

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.885806666666667
$ws.Range("H5").Value = 5.65742
$ws.Range("I5").Value = 0.3170188432052781
$ws.Range("J5").Value = 0.3170188432052781
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.15544366666667
$ws.Range("N5").Value = 75.466331
$ws.Range("O5").Value = 0.9701024243751556
$ws.Range("P5").Value = 0.9701024243751556
$ws.Range("Q5").Value = 47.43830336955778
$ws.Range("R5").Value = 426.94473032602
$ws.Range("S5").Value = 0.3075407483660476
$ws.Range("T5").Value = 0.3075407483660476

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.885806666666667
$ws.Range("H6").Value = 5.65742
$ws.Range("I6").Value = 0.3170188432052781
$ws.Range("J6").Value = 0.3170188432052781
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.6301496666666666
$ws.Range("N6").Value = 1.890449
$ws.Range("O6").Value = 0.02430128951224074
$ws.Range("P6").Value = 0.02430128951224074
$ws.Range("Q6").Value = 1.188340442397778
$ws.Range("R6").Value = 10.69506398158
$ws.Range("S6").Value = 0.007703966689567117
$ws.Range("T6").Value = 0.007703966689567117

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.885806666666667
$ws.Range("H7").Value = 5.65742
$ws.Range("I7").Value = 0.3170188432052781
$ws.Range("J7").Value = 0.3170188432052781
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.1451156666666667
$ws.Range("N7").Value = 0.435347
$ws.Range("O7").Value = 0.005596286112603657
$ws.Range("P7").Value = 0.005596286112603657
$ws.Range("Q7").Value = 0.2736600916377778
$ws.Range("R7").Value = 2.46294082474
$ws.Range("S7").Value = 0.001774128149663374
$ws.Range("T7").Value = 0.001774128149663374
